$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "C123413"
$ws.Range("B2").Value = "abcde"
